# Refactor lecturer and module timetable without redundant version
#
# This script rewrites a handful of cells on the timetable grid (moving /
# correcting several tutorial & practical sessions) and adjusts the
# corresponding merged cell ranges.
#
# Fill styles in the source workbook come from shared cellXfs entries
# (style "4" = tutorial-block fill, style "5" = practical-block fill).
# Rather than re-creating those fills by hand (which would mint brand-new
# style entries in the saved package and risk mismatched borders), we copy
# the *format only* from cells that already carry the exact style we need
# via PasteSpecial(xlPasteFormats). Cells that need to end up empty get
# their format copied from a blank, untouched scratch cell so they
# serialize back out with no style override at all (same as a cell that
# was never written).
#
# IMPORTANT ordering constraint discovered empirically against this
# engine's merge implementation: Range.Merge() redistributes borders
# across the newly merged block (the left-most cell loses its right
# border, the interior cells get a separate borderless style) whenever the
# destination cells already carry a "previous" style/value. To keep every
# merged block reporting as a single anchor cell with the plain, shared
# style (matching how the rest of the sheet's merges already look), all
# Merge()/UnMerge() calls are performed FIRST, while the affected cells are
# still in their original ("before") state, and only afterwards do we
# overwrite cell values/formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Reference cells whose style never changes - used as format donors.
$style4Source = "J2"   # tutorial-style fill (s="4")
$style5Source = "F2"   # practical-style fill (s="5")
$blankSource  = "Z1"   # untouched cell far outside the used range -> default/no style

# --- Step 1: fix up merged ranges first (before touching cell content) ---
$ws.Range("K4:N4").UnMerge()
$ws.Range("R5:U5").UnMerge()
$ws.Range("G10:J10").UnMerge()
$ws.Range("J11:M11").UnMerge()

$ws.Range("R3:U3").Merge()
$ws.Range("J4:M4").Merge()
$ws.Range("K5:N5").Merge()
$ws.Range("J10:M10").Merge()

# --- Step 2: now update the cell contents/formats ---
function Set-CellText($ref, $text) {
    $ws.Range($ref).Value2 = $text
}

function Copy-Format($fromRef, $toRef) {
    $ws.Range($fromRef).Copy() | Out-Null
    $ws.Range($toRef).PasteSpecial($xlPasteFormats) | Out-Null
}

function Set-TutorialCell($ref, $text) {
    Copy-Format $style4Source $ref
    Set-CellText $ref $text
}

function Set-PracticalCell($ref, $text) {
    Copy-Format $style5Source $ref
    Set-CellText $ref $text
}

function Clear-Cell($ref) {
    # Paste blank formatting first, then blank the value, so the cell
    # round-trips as if it was never written (no stray <c> element).
    Copy-Format $blankSource $ref
    Set-CellText $ref ""
}

# --- Monday (row 2-3): add a new MTH1114-G5-Tutorial slot at R3 ---
Set-TutorialCell "R3" "MTH1114-G5-Tutorial`n(Jaya Krishna)`nUE2-17"

# --- Tuesday (rows 4-5) ---
# New ENG1044-G3-Tutorial slot (correct room UE2-17) at J4.
Set-TutorialCell "J4" "ENG1044-G3-Tutorial`n(Herrick Yeap Han Lin)`nUE2-17"
# CSC1202-G3-Practical moves from K4 down to K5.
Set-PracticalCell "K5" "CSC1202-G3-Practical`n(Chew Moi Tin)`nUE2-16"
Clear-Cell "K4"
# F5 becomes MTH1114-G6-Tutorial; R5 (old MTH1114-G6 slot) is vacated.
Set-TutorialCell "F5" "MTH1114-G6-Tutorial`n(Jaya Krishna)`nUE2-17"
Clear-Cell "R5"

# --- Thursday (rows 7-8): ENG1044-G4-Tutorial and CSC1024-G3-Practical swap rows ---
Set-TutorialCell "N7" "ENG1044-G4-Tutorial`n(Herrick Yeap Han Lin)`nUC3-3"
Set-PracticalCell "N8" "CSC1024-G3-Practical`n(Farrukh Hassan)`nUE2-16"

# --- Friday (rows 9-11) ---
Set-TutorialCell "B9" "ENG1044-G5-Tutorial`n(Herrick Yeap Han Lin)`nUC3-3"
Set-PracticalCell "G9" "CSC1024-G4-Practical`n(Tan Kai Wei)`nUE2-16"
Set-TutorialCell "B10" "MTH1114-G1-Tutorial`n(Jaya Krishna)`nUE2-17"
Clear-Cell "G10"
Set-TutorialCell "J10" "MTH1114-G2-Tutorial`n(Jaya Krishna)`nUE2-17"
Clear-Cell "J11"
